# Apply crypto price/volume updates from Fri Oct 27 10:56:26 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.082.92'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '1.784.51'
$ws.Range("E3").Value = '  -2.36%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("E6").Value = '  -1.19%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.79'
$ws.Range("E8").Value = '  +2.45%  '

$ws.Range("E9").Value = '  -1.97%  '

$ws.Range("E10").Value = '  -0.72%  '

$ws.Range("E11").Value = '  +0.70%  '

$ws.Range("D12").Value = '2.041.79'
$ws.Range("E12").Value = '  -2.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.07'
$ws.Range("E13").Value = '  +2.23%  '

$ws.Range("D14").Value = '1.784.88'
$ws.Range("E14").Value = '  -2.66%  '

$ws.Range("D15").Value = '34.030.30'
$ws.Range("E15").Value = '  -0.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.620'
$ws.Range("E16").Value = '  -3.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.16'
$ws.Range("E17").Value = '  -4.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.91'
$ws.Range("E18").Value = '  -2.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.88'
$ws.Range("E19").Value = '  -2.74%  '

$ws.Range("D20").Value = '0.0₃0787'
$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.76'
$ws.Range("E22").Value = '  -4.01%  '

$ws.Range("E23").Value = '  -4.38%  '

$ws.Range("E24").Value = '  -2.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.49'
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("E26").Value = '  -1.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.05'
$ws.Range("E27").Value = '  -2.90%  '

$ws.Range("E28").Value = '  -2.09%  '

$ws.Range("E30").Value = '  +0.78%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0513'
$ws.Range("E31").Value = '  -4.29%  '

$ws.Range("E32").Value = '  -3.98%  '

$ws.Range("E33").Value = '  -1.83%  '

$ws.Range("E34").Value = '  -4.61%  '

$ws.Range("D35").Value = '1.390.15'
$ws.Range("E35").Value = '  -3.52%  '

$ws.Range("E36").Value = '  +0.23%  '

$ws.Range("E37").Value = '  -1.50%  '

$ws.Range("E38").Value = '  -1.74%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.20'
$ws.Range("E39").Value = '  +2.52%  '

$ws.Range("E40").Value = '  -0.20%  '

$ws.Range("E41").Value = '  -5.20%  '

$ws.Range("E42").Value = '  -2.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '78.18'
$ws.Range("E43").Value = '  -4.29%  '

$ws.Range("D44").Value = '0.0₆0142'
$ws.Range("E44").Value = '  +14.62%  '

$ws.Range("E45").Value = '  +2.27%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.52'
$ws.Range("E46").Value = '  +5.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0498'
$ws.Range("E47").Value = '  -0.15%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '108.04'
$ws.Range("E48").Value = '  +1.51%  '

$ws.Range("E49").Value = '  -4.12%  '

$ws.Range("D50").Value = '1.941.97'
$ws.Range("E50").Value = '  -2.36%  '

$ws.Range("E51").Value = '  +0.05%  '
